$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from existing header cell (AC1) to new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill data rows 2-47 with Wins=76, Losses=86, Ties=0
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 76  # AD
    $ws.Cells.Item($r, 31).Value = 86  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
